$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("G2").Value = 18.301715
$ws.Range("H2").Value = 36.60343
$ws.Range("I2").Value = 0.4544165459658234
$ws.Range("J2").Value = 0.3632302676703625
$ws.Range("K2").Value = 2
$ws.Range("M2").Value = 55.873922
$ws.Range("N2").Value = 111.747844
$ws.Range("O2").Value = 0.5825754174893317
$ws.Range("P2").Value = 0.4843347495294592
$ws.Range("Q2").Value = 1022.58859637623
$ws.Range("R2").Value = 4090.35438550492
$ws.Range("S2").Value = 0.2647319089800996
$ws.Range("T2").Value = 0.1759250407136434
$ws.Range("E3").Value = 2
$ws.Range("G3").Value = 18.301715
$ws.Range("H3").Value = 36.60343
$ws.Range("I3").Value = 0.4544165459658234
$ws.Range("J3").Value = 0.3632302676703625
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.743428333333333
$ws.Range("N3").Value = 5.230285
$ws.Range("O3").Value = 0.01817804179120264
$ws.Range("P3").Value = 0.02266897225724272
$ws.Range("Q3").Value = 31.90772847959167
$ws.Range("R3").Value = 191.44637087755
$ws.Range("S3").Value = 0.008260402963180691
$ws.Range("T3").Value = 0.008234056860810295
$ws.Range("E4").Value = 2
$ws.Range("G4").Value = 18.301715
$ws.Range("H4").Value = 36.60343
$ws.Range("I4").Value = 0.4544165459658234
$ws.Range("J4").Value = 0.3632302676703625
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.001563
$ws.Range("N4").Value = 3.004689
$ws.Range("O4").Value = 0.0104429036298341
$ws.Range("P4").Value = 0.01302284896189067
$ws.Range("Q4").Value = 18.330320580545
$ws.Range("R4").Value = 109.98192348327
$ws.Range("S4").Value = 0.004745428197323173
$ws.Range("T4").Value = 0.004730292914258252
$ws.Range("E5").Value = 2
$ws.Range("G5").Value = 18.301715
$ws.Range("H5").Value = 36.60343
$ws.Range("I5").Value = 0.4544165459658234
$ws.Range("J5").Value = 0.3632302676703625
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.365389333333333
$ws.Range("N5").Value = 4.096168
$ws.Range("O5").Value = 0.014236377766754
$ws.Range("P5").Value = 0.01775351032553778
$ws.Range("Q5").Value = 24.98896644270667
$ws.Range("R5").Value = 149.93379865624
$ws.Range("S5").Value = 0.006469245611832992
$ws.Range("T5").Value = 0.006448612307633633
$ws.Range("E6").Value = 2
$ws.Range("G6").Value = 18.301715
$ws.Range("H6").Value = 36.60343
$ws.Range("I6").Value = 0.4544165459658234
$ws.Range("J6").Value = 0.3632302676703625
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 34.797061
$ws.Range("N6").Value = 104.391183
$ws.Range("O6").Value = 0.3628152743506486
$ws.Range("P6").Value = 0.4524496908538919
$ws.Range("Q6").Value = 636.8458932596151
$ws.Range("R6").Value = 3821.07535955769
$ws.Range("S6").Value = 0.1648692637940644
$ws.Range("T6").Value = 0.1643434223162319
$ws.Range("E7").Value = 2
$ws.Range("G7").Value = 18.301715
$ws.Range("H7").Value = 36.60343
$ws.Range("I7").Value = 0.4544165459658234
$ws.Range("J7").Value = 0.3632302676703625
$ws.Range("K7").Value = 2
$ws.Range("M7").Value = 1.127115
$ws.Range("N7").Value = 2.25423
$ws.Range("O7").Value = 0.01175198497222887
$ws.Range("P7").Value = 0.00977022807197777
$ws.Range("Q7").Value = 20.628137502225
$ws.Range("R7").Value = 82.51255000889999
$ws.Range("S7").Value = 0.005340296419322506
$ws.Range("T7").Value = 0.003548842557784975
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.2426546666666667
$ws.Range("H8").Value = 0.7279639999999999
$ws.Range("I8").Value = 0.006024915997716873
$ws.Range("J8").Value = 0.007223873789270233
$ws.Range("K8").Value = 2
$ws.Range("M8").Value = 55.873922
$ws.Range("N8").Value = 111.747844
$ws.Range("O8").Value = 0.5825754174893317
$ws.Range("P8").Value = 0.4843347495294592
$ws.Range("Q8").Value = 13.55806791826933
$ws.Range("R8").Value = 81.348407509616
$ws.Range("S8").Value = 0.003509967952708061
$ws.Range("T8").Value = 0.003498773102358623
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.2426546666666667
$ws.Range("H9").Value = 0.7279639999999999
$ws.Range("I9").Value = 0.006024915997716873
$ws.Range("J9").Value = 0.007223873789270233
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.743428333333333
$ws.Range("N9").Value = 5.230285
$ws.Range("O9").Value = 0.01817804179120264
$ws.Range("P9").Value = 0.02266897225724272
$ws.Range("Q9").Value = 0.4230510210822222
$ws.Range("R9").Value = 3.80745918974
$ws.Range("S9").Value = 0.0001095211747949827
$ws.Range("T9").Value = 0.0001637577945187898
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.2426546666666667
$ws.Range("H10").Value = 0.7279639999999999
$ws.Range("I10").Value = 0.006024915997716873
$ws.Range("J10").Value = 0.007223873789270233
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.001563
$ws.Range("N10").Value = 3.004689
$ws.Range("O10").Value = 0.0104429036298341
$ws.Range("P10").Value = 0.01302284896189067
$ws.Range("Q10").Value = 0.2430339359106667
$ws.Range("R10").Value = 2.187305423196
$ws.Range("S10").Value = 0.00006291761714200309
$ws.Range("T10").Value = 0.00009407541727742711
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.2426546666666667
$ws.Range("H11").Value = 0.7279639999999999
$ws.Range("I11").Value = 0.006024915997716873
$ws.Range("J11").Value = 0.007223873789270233
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 1.365389333333333
$ws.Range("N11").Value = 4.096168
$ws.Range("O11").Value = 0.014236377766754
$ws.Range("P11").Value = 0.01775351032553778
$ws.Range("Q11").Value = 0.3313180935502222
$ws.Range("R11").Value = 2.981862841952
$ws.Range("S11").Value = 0.00008577298015645695
$ws.Range("T11").Value = 0.0001282491179081908
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.2426546666666667
$ws.Range("H12").Value = 0.7279639999999999
$ws.Range("I12").Value = 0.006024915997716873
$ws.Range("J12").Value = 0.007223873789270233
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 34.797061
$ws.Range("N12").Value = 104.391183
$ws.Range("O12").Value = 0.3628152743506486
$ws.Range("P12").Value = 0.4524496908538919
$ws.Range("Q12").Value = 8.443669237934666
$ws.Range("R12").Value = 75.993023141412
$ws.Range("S12").Value = 0.002185931550651259
$ws.Range("T12").Value = 0.003268439462722849
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.2426546666666667
$ws.Range("H13").Value = 0.7279639999999999
$ws.Range("I13").Value = 0.006024915997716873
$ws.Range("J13").Value = 0.007223873789270233
$ws.Range("K13").Value = 2
$ws.Range("M13").Value = 1.127115
$ws.Range("N13").Value = 2.25423
$ws.Range("O13").Value = 0.01175198497222887
$ws.Range("P13").Value = 0.00977022807197777
$ws.Range("Q13").Value = 0.2734997146199999
$ws.Range("R13").Value = 1.64099828772
$ws.Range("S13").Value = 0.00007080472226411001
$ws.Range("T13").Value = 0.00007057889448435246
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 1.82251
$ws.Range("H14").Value = 5.46753
$ws.Range("I14").Value = 0.0452514258466036
$ws.Range("J14").Value = 0.05425645589486387
$ws.Range("K14").Value = 2
$ws.Range("M14").Value = 55.873922
$ws.Range("N14").Value = 111.747844
$ws.Range("O14").Value = 0.5825754174893317
$ws.Range("P14").Value = 0.4843347495294592
$ws.Range("Q14").Value = 101.83078158422
$ws.Range("R14").Value = 610.98468950532
$ws.Range("S14").Value = 0.02636236830457263
$ws.Range("T14").Value = 0.02627828697619504
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 1.82251
$ws.Range("H15").Value = 5.46753
$ws.Range("I15").Value = 0.0452514258466036
$ws.Range("J15").Value = 0.05425645589486387
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 1.743428333333333
$ws.Range("N15").Value = 5.230285
$ws.Range("O15").Value = 0.01817804179120264
$ws.Range("P15").Value = 0.02266897225724272
$ws.Range("Q15").Value = 3.177415571783333
$ws.Range("R15").Value = 28.59674014605
$ws.Range("S15").Value = 0.0008225823101510674
$ws.Range("T15").Value = 0.001229938093456982
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 1.82251
$ws.Range("H16").Value = 5.46753
$ws.Range("I16").Value = 0.0452514258466036
$ws.Range("J16").Value = 0.05425645589486387
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 1.001563
$ws.Range("N16").Value = 3.004689
$ws.Range("O16").Value = 0.0104429036298341
$ws.Range("P16").Value = 0.01302284896189067
$ws.Range("Q16").Value = 1.82535858313
$ws.Range("R16").Value = 16.42822724817
$ws.Range("S16").Value = 0.0004725562792286655
$ws.Range("T16").Value = 0.0007065736303262951
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 1.82251
$ws.Range("H17").Value = 5.46753
$ws.Range("I17").Value = 0.0452514258466036
$ws.Range("J17").Value = 0.05425645589486387
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 1.365389333333333
$ws.Range("N17").Value = 4.096168
$ws.Range("O17").Value = 0.014236377766754
$ws.Range("P17").Value = 0.01775351032553778
$ws.Range("Q17").Value = 2.488435713893333
$ws.Range("R17").Value = 22.39592142504
$ws.Range("S17").Value = 0.0006442163928365045
$ws.Range("T17").Value = 0.000963242549956551
$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 1.82251
$ws.Range("H18").Value = 5.46753
$ws.Range("I18").Value = 0.0452514258466036
$ws.Range("J18").Value = 0.05425645589486387
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 34.797061
$ws.Range("N18").Value = 104.391183
$ws.Range("O18").Value = 0.3628152743506486
$ws.Range("P18").Value = 0.4524496908538919
$ws.Range("Q18").Value = 63.41799164311
$ws.Range("R18").Value = 570.76192478799
$ws.Range("S18").Value = 0.01641790848329352
$ws.Range("T18").Value = 0.02454831669645897
$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 1.82251
$ws.Range("H19").Value = 5.46753
$ws.Range("I19").Value = 0.0452514258466036
$ws.Range("J19").Value = 0.05425645589486387
$ws.Range("K19").Value = 2
$ws.Range("M19").Value = 1.127115
$ws.Range("N19").Value = 2.25423
$ws.Range("O19").Value = 0.01175198497222887
$ws.Range("P19").Value = 0.00977022807197777
$ws.Range("Q19").Value = 2.05417835865
$ws.Range("R19").Value = 12.3250701519
$ws.Range("S19").Value = 0.0005317940765212145
$ws.Range("T19").Value = 0.0005300979484700227
$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 4.296528666666666
$ws.Range("H20").Value = 12.889586
$ws.Range("I20").Value = 0.1066792765787147
$ws.Range("J20").Value = 0.1279084439064906
$ws.Range("K20").Value = 2
$ws.Range("M20").Value = 55.873922
$ws.Range("N20").Value = 111.747844
$ws.Range("O20").Value = 0.5825754174893317
$ws.Range("P20").Value = 0.4843347495294592
$ws.Range("Q20").Value = 240.0639075920973
$ws.Range("R20").Value = 1440.383445552584
$ws.Range("S20").Value = 0.06214872409030459
$ws.Range("T20").Value = 0.06195050414215301
$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 4.296528666666666
$ws.Range("H21").Value = 12.889586
$ws.Range("I21").Value = 0.1066792765787147
$ws.Range("J21").Value = 0.1279084439064906
$ws.Range("K21").Value = 3
$ws.Range("M21").Value = 1.743428333333333
$ws.Range("N21").Value = 5.230285
$ws.Range("O21").Value = 0.01817804179120264
$ws.Range("P21").Value = 0.02266897225724272
$ws.Range("Q21").Value = 7.490689812445555
$ws.Range("R21").Value = 67.41620831201
$ws.Range("S21").Value = 0.00193922034790314
$ws.Range("T21").Value = 0.002899552966383322
$ws.Range("E22").Value = 3
$ws.Range("G22").Value = 4.296528666666666
$ws.Range("H22").Value = 12.889586
$ws.Range("I22").Value = 0.1066792765787147
$ws.Range("J22").Value = 0.1279084439064906
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 1.001563
$ws.Range("N22").Value = 3.004689
$ws.Range("O22").Value = 0.0104429036298341
$ws.Range("P22").Value = 0.01302284896189067
$ws.Range("Q22").Value = 4.303244140972666
$ws.Range("R22").Value = 38.729197268754
$ws.Range("S22").Value = 0.001114041404611936
$ws.Range("T22").Value = 0.001665732345944693
$ws.Range("E23").Value = 3
$ws.Range("G23").Value = 4.296528666666666
$ws.Range("H23").Value = 12.889586
$ws.Range("I23").Value = 0.1066792765787147
$ws.Range("J23").Value = 0.1279084439064906
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 1.365389333333333
$ws.Range("N23").Value = 4.096168
$ws.Range("O23").Value = 0.014236377766754
$ws.Range("P23").Value = 0.01775351032553778
$ws.Range("Q23").Value = 5.866434411827554
$ws.Range("R23").Value = 52.797909706448
$ws.Range("S23").Value = 0.001518726481258614
$ws.Range("T23").Value = 0.002270823879617352
$ws.Range("E24").Value = 3
$ws.Range("G24").Value = 4.296528666666666
$ws.Range("H24").Value = 12.889586
$ws.Range("I24").Value = 0.1066792765787147
$ws.Range("J24").Value = 0.1279084439064906
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 34.797061
$ws.Range("N24").Value = 104.391183
$ws.Range("O24").Value = 0.3628152743506486
$ws.Range("P24").Value = 0.4524496908538919
$ws.Range("Q24").Value = 149.5065701022486
$ws.Range("R24").Value = 1345.559130920238
$ws.Range("S24").Value = 0.03870487099943509
$ws.Range("T24").Value = 0.05787213590309405
$ws.Range("E25").Value = 3
$ws.Range("G25").Value = 4.296528666666666
$ws.Range("H25").Value = 12.889586
$ws.Range("I25").Value = 0.1066792765787147
$ws.Range("J25").Value = 0.1279084439064906
$ws.Range("K25").Value = 2
$ws.Range("M25").Value = 1.127115
$ws.Range("N25").Value = 2.25423
$ws.Range("O25").Value = 0.01175198497222887
$ws.Range("P25").Value = 0.00977022807197777
$ws.Range("Q25").Value = 4.842681908129999
$ws.Range("R25").Value = 29.05609144877999
$ws.Range("S25").Value = 0.001253693255201302
$ws.Range("T25").Value = 0.001249694669298188
$ws.Range("E26").Value = 3
$ws.Range("G26").Value = 13.85988766666667
$ws.Range("H26").Value = 41.579663
$ws.Range("I26").Value = 0.3441296228774725
$ws.Range("J26").Value = 0.4126113897285982
$ws.Range("K26").Value = 2
$ws.Range("M26").Value = 55.873922
$ws.Range("N26").Value = 111.747844
$ws.Range("O26").Value = 0.5825754174893317
$ws.Range("P26").Value = 0.4843347495294592
$ws.Range("Q26").Value = 774.4062824160952
$ws.Range("R26").Value = 4646.437694496572
$ws.Range("S26").Value = 0.2004814587182898
$ws.Range("T26").Value = 0.1998420340972027
$ws.Range("E27").Value = 3
$ws.Range("G27").Value = 13.85988766666667
$ws.Range("H27").Value = 41.579663
$ws.Range("I27").Value = 0.3441296228774725
$ws.Range("J27").Value = 0.4126113897285982
$ws.Range("K27").Value = 3
$ws.Range("M27").Value = 1.743428333333333
$ws.Range("N27").Value = 5.230285
$ws.Range("O27").Value = 0.01817804179120264
$ws.Range("P27").Value = 0.02266897225724272
$ws.Range("Q27").Value = 24.16372085488389
$ws.Range("R27").Value = 217.473487693955
$ws.Range("S27").Value = 0.006255602666257497
$ws.Range("T27").Value = 0.009353476146779957
$ws.Range("E28").Value = 3
$ws.Range("G28").Value = 13.85988766666667
$ws.Range("H28").Value = 41.579663
$ws.Range("I28").Value = 0.3441296228774725
$ws.Range("J28").Value = 0.4126113897285982
$ws.Range("K28").Value = 3
$ws.Range("M28").Value = 1.001563
$ws.Range("N28").Value = 3.004689
$ws.Range("O28").Value = 0.0104429036298341
$ws.Range("P28").Value = 0.01302284896189067
$ws.Range("Q28").Value = 13.88155067108966
$ws.Range("R28").Value = 124.933956039807
$ws.Range("S28").Value = 0.003593712487880598
$ws.Range("T28").Value = 0.005373375808391344
$ws.Range("E29").Value = 3
$ws.Range("G29").Value = 13.85988766666667
$ws.Range("H29").Value = 41.579663
$ws.Range("I29").Value = 0.3441296228774725
$ws.Range("J29").Value = 0.4126113897285982
$ws.Range("K29").Value = 3
$ws.Range("M29").Value = 1.365389333333333
$ws.Range("N29").Value = 4.096168
$ws.Range("O29").Value = 0.014236377766754
$ws.Range("P29").Value = 0.01775351032553778
$ws.Range("Q29").Value = 18.92414278126488
$ws.Range("R29").Value = 170.317285031384
$ws.Range("S29").Value = 0.004899159312014286
$ws.Range("T29").Value = 0.007325300567981164
$ws.Range("E30").Value = 3
$ws.Range("G30").Value = 13.85988766666667
$ws.Range("H30").Value = 41.579663
$ws.Range("I30").Value = 0.3441296228774725
$ws.Range("J30").Value = 0.4126113897285982
$ws.Range("K30").Value = 3
$ws.Range("M30").Value = 34.797061
$ws.Range("N30").Value = 104.391183
$ws.Range("O30").Value = 0.3628152743506486
$ws.Range("P30").Value = 0.4524496908538919
$ws.Range("Q30").Value = 482.2833565901476
$ws.Range("R30").Value = 4340.550209311328
$ws.Range("S30").Value = 0.1248554835364754
$ws.Range("T30").Value = 0.186685895725499
$ws.Range("E31").Value = 3
$ws.Range("G31").Value = 13.85988766666667
$ws.Range("H31").Value = 41.579663
$ws.Range("I31").Value = 0.3441296228774725
$ws.Range("J31").Value = 0.4126113897285982
$ws.Range("K31").Value = 2
$ws.Range("M31").Value = 1.127115
$ws.Range("N31").Value = 2.25423
$ws.Range("O31").Value = 0.01175198497222887
$ws.Range("P31").Value = 0.00977022807197777
$ws.Range("Q31").Value = 15.621687287415
$ws.Range("R31").Value = 93.73012372448999
$ws.Range("S31").Value = 0.004044206156554844
$ws.Range("T31").Value = 0.00403130738274411
$ws.Range("E32").Value = 2
$ws.Range("G32").Value = 1.751899
$ws.Range("H32").Value = 3.503798
$ws.Range("I32").Value = 0.04349821273366894
$ws.Range("J32").Value = 0.03476956901041461
$ws.Range("K32").Value = 2
$ws.Range("M32").Value = 55.873922
$ws.Range("N32").Value = 111.747844
$ws.Range("O32").Value = 0.5825754174893317
$ws.Range("P32").Value = 0.4843347495294592
$ws.Range("Q32").Value = 97.88546807787799
$ws.Range("R32").Value = 391.541872311512
$ws.Range("S32").Value = 0.02534098944335695
$ws.Range("T32").Value = 0.01684011049790641
$ws.Range("E33").Value = 2
$ws.Range("G33").Value = 1.751899
$ws.Range("H33").Value = 3.503798
$ws.Range("I33").Value = 0.04349821273366894
$ws.Range("J33").Value = 0.03476956901041461
$ws.Range("K33").Value = 3
$ws.Range("M33").Value = 1.743428333333333
$ws.Range("N33").Value = 5.230285
$ws.Range("O33").Value = 0.01817804179120264
$ws.Range("P33").Value = 0.02266897225724272
$ws.Range("Q33").Value = 3.054310353738333
$ws.Range("R33").Value = 18.32586212243
$ws.Range("S33").Value = 0.0007907123289152566
$ws.Range("T33").Value = 0.0007881903952933751
$ws.Range("E34").Value = 2
$ws.Range("G34").Value = 1.751899
$ws.Range("H34").Value = 3.503798
$ws.Range("I34").Value = 0.04349821273366894
$ws.Range("J34").Value = 0.03476956901041461
$ws.Range("K34").Value = 3
$ws.Range("M34").Value = 1.001563
$ws.Range("N34").Value = 3.004689
$ws.Range("O34").Value = 0.0104429036298341
$ws.Range("P34").Value = 0.01302284896189067
$ws.Range("Q34").Value = 1.754637218137
$ws.Range("R34").Value = 10.527823308822
$ws.Range("S34").Value = 0.0004542476436477274
$ws.Range("T34").Value = 0.0004527988456926641
$ws.Range("E35").Value = 2
$ws.Range("G35").Value = 1.751899
$ws.Range("H35").Value = 3.503798
$ws.Range("I35").Value = 0.04349821273366894
$ws.Range("J35").Value = 0.03476956901041461
$ws.Range("K35").Value = 3
$ws.Range("M35").Value = 1.365389333333333
$ws.Range("N35").Value = 4.096168
$ws.Range("O35").Value = 0.014236377766754
$ws.Range("P35").Value = 0.01775351032553778
$ws.Range("Q35").Value = 2.392024207677333
$ws.Range("R35").Value = 14.352145246064
$ws.Range("S35").Value = 0.00061925698865514
$ws.Range("T35").Value = 0.0006172819024408944
$ws.Range("E36").Value = 2
$ws.Range("G36").Value = 1.751899
$ws.Range("H36").Value = 3.503798
$ws.Range("I36").Value = 0.04349821273366894
$ws.Range("J36").Value = 0.03476956901041461
$ws.Range("K36").Value = 3
$ws.Range("M36").Value = 34.797061
$ws.Range("N36").Value = 104.391183
$ws.Range("O36").Value = 0.3628152743506486
$ws.Range("P36").Value = 0.4524496908538919
$ws.Range("Q36").Value = 60.96093636883899
$ws.Range("R36").Value = 365.765618213034
$ws.Range("S36").Value = 0.01578181598672897
$ws.Range("T36").Value = 0.01573148074988515
$ws.Range("E37").Value = 2
$ws.Range("G37").Value = 1.751899
$ws.Range("H37").Value = 3.503798
$ws.Range("I37").Value = 0.04349821273366894
$ws.Range("J37").Value = 0.03476956901041461
$ws.Range("K37").Value = 2
$ws.Range("M37").Value = 1.127115
$ws.Range("N37").Value = 2.25423
$ws.Range("O37").Value = 0.01175198497222887
$ws.Range("P37").Value = 0.00977022807197777
$ws.Range("Q37").Value = 1.974591641385
$ws.Range("R37").Value = 7.898366565539998
$ws.Range("S37").Value = 0.0005111903423648919
$ws.Range("T37").Value = 0.0003397066191961212
